$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore it at the end
# (row inserts / selection changes below must not change the active tab).
$origActiveSheetName = $wb.ActiveSheet.Name

# --- Sheet "p2" (sheet3.xml): insert a new blank row above row 5 ---
$wsP2 = $wb.Worksheets.Item("p2")
$wsP2.Rows.Item(5).Insert() | Out-Null
$wsP2.Rows.Item(5).RowHeight = 18.75
$wsP2.Activate() | Out-Null
$wsP2.Range("B11").Select() | Out-Null

# --- Sheet "p3" (sheet4.xml): insert a new blank row above row 5 ---
$wsP3 = $wb.Worksheets.Item("p3")
$wsP3.Rows.Item(5).Insert() | Out-Null
$wsP3.Rows.Item(5).RowHeight = 18.75
$wsP3.Activate() | Out-Null
$wsP3.Range("B13").Select() | Out-Null

# --- Sheet "p4" (sheet5.xml): insert a new blank row above row 5 ---
$wsP4 = $wb.Worksheets.Item("p4")
$wsP4.Rows.Item(5).Insert() | Out-Null
$wsP4.Rows.Item(5).RowHeight = 18.75
$wsP4.Activate() | Out-Null
$wsP4.Range("B11").Select() | Out-Null

# Restore the originally active sheet/tab
$wb.Worksheets.Item($origActiveSheetName).Activate() | Out-Null
